$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUP_TRADE")

# --- Update input prices (column H, first trade block) ---
$ws.Range("H5").Value = 0.0895
$ws.Range("H6").Value = 0.00006
$ws.Range("H9").Value = 20
$ws.Range("H10").Value = 0.0031329

# H11 changes value AND gets a new scientific number format (0.00E+00),
# which also re-styles the dependent cells that reference it (I21/I31/I41)
$ws.Range("H11").Value = 0.00001566
$ws.Range("H11").NumberFormat = "0.00E+00"

# --- Second trade block (column I) ---
$ws.Range("I15").Value = 0.2983
$ws.Range("I16").Value = 0.00019
$ws.Range("I19").Value = 20

# I20 becomes a formula referencing H10 instead of a literal value
$ws.Range("I20").Formula = "=H10"

# I21 becomes a formula referencing H11, with the same new number format
$ws.Range("I21").Formula = "=H11"
$ws.Range("I21").NumberFormat = "0.00E+00"

# --- Third trade block (column I) ---
$ws.Range("I25").Value = 0.3878
$ws.Range("I26").Value = 0.00024
$ws.Range("I29").Value = 20

# I30 becomes a formula referencing I20
$ws.Range("I30").Formula = "=I20"

# I31 becomes a formula referencing I21, with the same new number format
$ws.Range("I31").Formula = "=I21"
$ws.Range("I31").NumberFormat = "0.00E+00"

# --- Fourth trade block (column I) ---
$ws.Range("I35").Value = 5370
$ws.Range("I36").Value = 0.00034
$ws.Range("I39").Value = 20

# I40 becomes a formula referencing I20
$ws.Range("I40").Formula = "=I20"

# I41 becomes a formula referencing I21, with the same new number format
$ws.Range("I41").Formula = "=I21"
$ws.Range("I41").NumberFormat = "0.00E+00"

# --- Update selected cell ---
$ws.Range("I42").Select() | Out-Null
